$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C1: replace the literal "\n" markers with real line breaks (Alt+Enter),
# keeping a real newline (char 10) between each segment, no trailing newline.
$newLine = [char]10
$newText = "Max. discount 200rb" + $newLine + "Without Min. Transactions" + $newLine + "Takeaway Only" + $newLine + "Valid until 15 May 2023"
$ws.Range("C1").Value = $newText

# Move the active selection to C1 (matches the saved selection state in the file)
$ws.Range("C1").Select()

$wb.Save()
